$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Bosses get stuck on ground when dashing" (row 5) is being re-triaged out of the
# "Highest Priority" section and moved down to the bottom of the list (uncategorised,
# newest entry) as part of adding music/audio related work items.
# Deleting the row shifts everything below it up by one.
$ws.Rows.Item(5).Delete() | Out-Null

# The blank spacer row that used to separate "Highest Priority" from the next section
# collapses as well (one less blank row is kept), shifting everything below it up by one more.
$ws.Rows.Item(3).Delete() | Out-Null

# Re-add "Bosses get stuck on ground when dashing" as a new, plain (unstyled) entry at
# the very end of the table.
$ws.Range("B94").Value = "Bosses get stuck on ground when dashing"

# Add a (currently blank-ish) note next to "Sound effects" in the Tutorials/Notes list.
$ws.Range("C10").Value = "  "

# Update the view: select C10, with no particular scroll position.
$ws.Activate() | Out-Null
$ws.Range("C10").Select() | Out-Null
